$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 534, pushing the existing row 534 (and everything
# below it) down by one. This mirrors the xlsx diff, where a new weekly
# price record is prepended ahead of the row that used to be first (old
# row 534 -> new row 535, ..., old row 641 -> new row 642).
$ws.Rows("534:534").Insert()

# Populate the newly inserted row 534 with the new record. Most fields
# mirror the (now-shifted) row below it (row 535, formerly row 534);
# only the date (D) and volume (J) differ, per the diff.
$ws.Range("A534").Value = 8
$ws.Range("B534").Value = "Terminal La Palmera de La Serena"
$ws.Range("C534").Value = "Coquimbo"
$ws.Range("D534").Value = 45258
$ws.Range("E534").Value = 4
$ws.Range("F534").Value = 100114013
$ws.Range("G534").Value = "Zanahoria"
$ws.Range("H534").Value = "Sin especificar"
$ws.Range("I534").Value = "Primera"
$ws.Range("J534").Value = 540
$ws.Range("K534").Value = 5500
$ws.Range("L534").Value = 6000
$ws.Range("M534").Value = 5750
$ws.Range("N534").Value = "$/saco 20 kilos"
$ws.Range("O534").Value = "Provincia del Elquí"
$ws.Range("P534").Value = 288
$ws.Range("Q534").Value = 20
$ws.Range("R534").Value = "Hortaliza"
